# Updated symbol list with refreshed Price (D) and Volume(1h) (E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value (kept as text to match existing inline-string cells)
$updates = @{
    "E2" = "-2.64%"
    "D3" = "31.07"
    "E3" = "-3.20%"
    "D4" = "4.938"
    "E4" = "-1.74%"
    "D5" = "0.07326"
    "E5" = "-7.12%"
    "D6" = "1.820"
    "E6" = "-15.25%"
    "D7" = "7.669"
    "E7" = "-1.86%"
    "D8" = "3.761"
    "E8" = "-0.94%"
    "D9" = "0.9067"
    "E9" = "-2.22%"
    "D10" = "0.1656"
    "E10" = "-5.24%"
    "D11" = "0.07559"
    "E11" = "-5.30%"
    "D12" = "0.08154"
    "E12" = "-7.09%"
    "D13" = "0.02988"
    "E13" = "-4.29%"
    "D14" = "0.09997"
    "E14" = "-0.42%"
    "D15" = "0.001491"
    "E15" = "-1.40%"
    "D16" = "0.005704"
    "E16" = "-2.06%"
    "E17" = "-0.17%"
    "D18" = "2.101"
    "E18" = "-7.79%"
    "D19" = "0.3272"
    "D20" = "0.1306"
    "E20" = "1.41%"
    "D21" = "4.374"
    "E21" = "5.73%"
    "D22" = "0.2000"
    "E22" = "11.87%"
    "D23" = "0.04478"
    "E23" = "-2.25%"
    "E24" = "-0.62%"
    "D25" = "0.004042"
    "E25" = "-10.62%"
    "E26" = "0.35%"
    "D39" = "0.01654"
    "E39" = "-5.40%"
    "D40" = "0.04395"
    "E40" = "-7.59%"
    "D41" = "0.007406"
    "E41" = "0.30%"
    "D42" = "0.1320"
    "E42" = "-3.61%"
    "D43" = "0.002091"
    "E43" = "-10.37%"
    "D44" = "0.01114"
    "E44" = "2.60%"
    "D45" = "0.00005991"
    "E45" = "-0.87%"
    "D46" = "0.00000000750"
    "E46" = "0.34%"
    "D47" = "2.138"
    "E47" = "159.62%"
    "D48" = "0.002401"
    "E48" = "-29.17%"
    "D49" = "0.00002101"
    "E49" = "0.34%"
    "D50" = "0.0002001"
    "E50" = "0.34%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text formatting so values such as "31.07" or "-2.64%" are stored
    # as plain text (matching the workbook original), not auto-converted to
    # numbers/percentages by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Revert the cell style back to the default "Normal" style so no stray
    # number-format styling remains applied to the cell itself.
    $cell.Style = "Normal"
}
